$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 45, shifting the existing rows 45-62 down to 46-63
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the "Il giro del mondo in 90 minuti" event
$ws.Range("A45").Value() = "Altri eventi,Spettacoli"
$ws.Range("B45").Value() = "Modena"
$ws.Range("C45").Value() = "ingresso da Strada S.Faustino, 172"
$ws.Range("D45").Value() = "2022-05-26T14:09:48+00:00"
$ws.Range("E45").Value() = "Spettacolo di danza"
$ws.Range("F45").Value() = "2022-05-26T14:09:58+00:00"
$ws.Range("G45").Value() = ""
$ws.Range("H45").Value() = "2022-06-04T14:00:00+00:00"
$ws.Range("I45").Value() = "2022-06-04T15:00:00+00:00"
$ws.Range("J45").Value() = "https://www.comune.modena.it/api/novita/eventi/2022/il-giro-del-mondo-in-90-minuti/@@images/c6bc9831-404e-402b-9e60-623bc0a2ba80.jpeg"
$ws.Range("K45").Value() = ""
$ws.Range("L45").Value() = "2022-05-26T14:11:15+00:00"
$ws.Range("M45").Value() = "Parco Ferrari"
$ws.Range("N45").Value() = " ore 20.00"
$ws.Range("O45").Value() = ""
$ws.Range("P45").Value() = " ingresso libero"
$ws.Range("Q45").Value() = ""
$ws.Range("R45").Value() = ""
$ws.Range("S45").Value() = "Il giro del mondo in 90 minuti"
$ws.Range("T45").Value() = ""
$ws.Range("U45").Value() = ""
$ws.Range("V45").Value() = $false
$ws.Range("W45").Value() = 41123
$ws.Range("X45").Value() = "https://www.comune.modena.it/novita/eventi/2022/il-giro-del-mondo-in-90-minuti"
$ws.Range("Y45").Value() = "44,64582"
$ws.Range("Z45").Value() = "10,92572"
$ws.Range("AA45").Value() = "POINT (10.92572 44.64582)"
